# Daily attendance processing - 2025-11-07 10:49:39
# Reorders the "Recorded By" (column G) author lists for rows where the
# list contains specific known value combinations, matching the upstream
# data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "System, admin@admin.com"              = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
